$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The "last_edited_time" value "2024-08-24T20:33:00.000Z" is stored once in the
# shared-string table and reused by every row that was last touched at that
# moment (D3, D4, D5, D7, D10, D13). The source edit updates that shared string
# in place, so all of these cells must be moved to the new timestamp together
# (otherwise Excel would keep the old string alive for the untouched cells and
# allocate a brand new shared string just for D10).
$newTimestamp = "2024-08-26T17:26:00.000Z"
$ws.Range("D3").Value = $newTimestamp
$ws.Range("D4").Value = $newTimestamp
$ws.Range("D5").Value = $newTimestamp
$ws.Range("D7").Value = $newTimestamp
$ws.Range("D10").Value = $newTimestamp
$ws.Range("D13").Value = $newTimestamp

# Update numeric figures for row 10
$ws.Range("T10").Value = 50500000
$ws.Range("W10").Value = 29618000
$ws.Range("AA10").Value = 83382000
$ws.Range("AK10").Value = 13
$ws.Range("AQ10").Value = 141500000
